$wb = $excel.ActiveWorkbook

# Rename sheets: H1a -> PL_H1a, H1b -> PL_H1b, H2b -> PL_H2b
$wb.Worksheets("H1a").Name = "PL_H1a"
$wb.Worksheets("H1b").Name = "PL_H1b"
$wb.Worksheets("H2b").Name = "PL_H2b"

# Activate PL_H2b (was H2b) -> becomes the selected/active tab (activeTab=3, tabSelected on that sheet)
$wb.Worksheets("PL_H2b").Activate()

# Update the active window size/position to match the new workbookView
$win = $excel.ActiveWindow
$win.Left = -120
$win.Top = -120
$win.Width = 38640
$win.Height = 21120
